$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.795.19"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "1.756.70"
$ws.Range("E3").Value = "  +0.81%  "

$ws.Range("E4").Value = "  +0.06%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.28"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("E6").Value = "  +0.02%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5066"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +3.27%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.25"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -1.62%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2619"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +8.48%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06187"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +2.82%  "

$ws.Range("D11").Value = "1.753.76"

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06949"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +4.61%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.44"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +6.65%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6008"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +1.28%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.60"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +0.30%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.446"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +2.57%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "25.841.00"
$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.63"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +3.47%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006800"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +8.00%  "

$ws.Range("D22").Value = "1.978.40"
$ws.Range("E22").Value = "  +0.00%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.056"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +4.58%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.136"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +2.35%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.165"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +1.45%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.85"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +1.81%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.460"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -4.21%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.00"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +5.39%  "

$ws.Range("E29").Value = "  -2.66%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.57"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +3.64%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08243"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.61%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.682"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +1.34%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.388"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +6.46%  "

$ws.Range("E34").Value = "  +1.74%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("D35").Style = $origStyle

$ws.Range("E36").Value = "  +1.68%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6007"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -0.82%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.729"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -2.58%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.936"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -6.48%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01545"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +3.91%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.36"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +2.37%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3811"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.09%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7435"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -5.90%  "

$ws.Range("E46").Value = "  -5.10%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05476"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +8.08%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1075"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +3.90%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.940"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -2.10%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.15"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +1.94%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9998"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +0.15%  "
